# Automatische test-sync: 2025-08-04 20:49:50
# Appends a new logged response row (row 7) to the historical responses sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Testmail #11: Mijn retour is nog steeds niet verwerkt."
$ws.Range("B7").Value = "Beste klant,`nDank u voor uw bericht. Om uw retourzending verder te kunnen onderzoeken, heb ik wat meer informatie nodig. Kunt u alstublieft uw ordernummer en de datum van de retourzending doorgeven? Op die manier kunnen wij u sneller van dienst zijn.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$ws.Range("C7").Value = "Mijn retour is nog steeds niet verwerkt."
$ws.Range("D7").Value = "mailmind.test@zohomail.eu"
$ws.Range("E7").Value = "Retour / Terugbetaling"
$ws.Range("F7").Value = "2025-08-04 20:49:17"
$ws.Range("G7").Value = "Ja"
$ws.Range("H7").Value = "Nee"
$ws.Range("I7").Value = "Ja"
$ws.Range("J7").Value = "Nee"

# The multi-line reply text in B7 causes Excel to auto-apply a custom row
# height; AutoFit() recomputes (and clears the manual-override flag) so the
# row is left at the sheet's default height, matching the other data rows.
$ws.Rows.Item(7).EntireRow.AutoFit()
